$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("E2").Value = 1

# Update row 6 values
$ws.Range("A6").Value = 1
$ws.Range("E6").Value = 0
